# This script updates the cryptocurrency price/volume table to reflect
# refreshed market data, including a re-ranking swap of several coins.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "67.953.88"
Set-TextCell "E2" "  +1.04%  "

Set-TextCell "D3" "2.540.46"
Set-TextCell "E3" "  +0.42%  "

Set-TextCell "E4" "  +0.01%  "

Set-TextCell "D5" "592.22"
Set-TextCell "E5" "  +0.24%  "

Set-TextCell "D6" "173.33"
Set-TextCell "E6" "  -0.20%  "

Set-TextCell "E7" "  -0.05%  "

Set-TextCell "E8" "  -0.32%  "

Set-TextCell "D9" "2.540.30"
Set-TextCell "E9" "  +0.45%  "

Set-TextCell "E10" "  +0.42%  "

Set-TextCell "E11" "  +1.77%  "

Set-TextCell "D12" "5.07"
Set-TextCell "E12" "  -1.63%  "

Set-TextCell "E13" "  -0.17%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D14" "26.48"
Set-TextCell "E14" "  -0.35%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell "D15" "2.973.66"
Set-TextCell "E15" "  -0.51%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D16" "0.0000178"
Set-TextCell "E16" "  +0.84%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D17" "67.856.25"
Set-TextCell "E17" "  +1.21%  "

$ws.Range("B18").Value = "Binance-PegBSC-USD"
$ws.Range("C18").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell "D18" "2.36"
Set-TextCell "E18" "  +136.54%  "

Set-TextCell "D19" "2.526.82"
Set-TextCell "E19" "  -0.58%  "

Set-TextCell "D20" "11.75"
Set-TextCell "E20" "  +3.13%  "

Set-TextCell "D21" "7.97"
Set-TextCell "E21" "  -1.35%  "

Set-TextCell "D22" "368.93"
Set-TextCell "E22" "  +4.00%  "

Set-TextCell "D23" "4.15"
Set-TextCell "E23" "  -0.72%  "

Set-TextCell "D24" "4.58"
Set-TextCell "E24" "  -0.82%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D25" "71.56"
Set-TextCell "E25" "  +2.51%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D26" "1.00"
Set-TextCell "E26" "  +0.10%  "

Set-TextCell "E27" "  -3.65%  "

Set-TextCell "D28" "9.95"
Set-TextCell "E28" "  -0.02%  "

Set-TextCell "E30" "  -0.97%  "

Set-TextCell "D31" "8.47"
Set-TextCell "E31" "  +3.66%  "

Set-TextCell "D32" "541.46"
Set-TextCell "E32" "  +1.52%  "

Set-TextCell "E33" "  -0.64%  "

Set-TextCell "E34" "  +1.08%  "

Set-TextCell "D35" "0.129"
Set-TextCell "E35" "  -1.03%  "

Set-TextCell "E36" "  -0.03%  "

Set-TextCell "D37" "159.29"
Set-TextCell "E37" "  +1.15%  "

Set-TextCell "E38" "  -1.90%  "

Set-TextCell "D39" "19.14"
Set-TextCell "E39" "  +2.59%  "

Set-TextCell "D40" "18.61"
Set-TextCell "E40" "  +0.87%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextCell "D41" "5.17"
Set-TextCell "E41" "  +0.56%  "

$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextCell "D42" "0.352"
Set-TextCell "E42" "  -0.95%  "

Set-TextCell "E43" "  -0.71%  "

Set-TextCell "E44" "  +2.82%  "

Set-TextCell "E45" "  -0.09%  "

Set-TextCell "D46" "39.29"
Set-TextCell "E46" "  -1.05%  "

Set-TextCell "D47" "0.0₆0288"
Set-TextCell "E47" "  +3.78%  "

Set-TextCell "E48" "  -0.93%  "

Set-TextCell "D49" "3.72"

Set-TextCell "E50" "  -0.82%  "

Set-TextCell "D51" "1.71"
Set-TextCell "E51" "  +1.06%  "

